$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 816.4
$ws.Range("J2").Value = 770.75
$ws.Range("L2").Value = 770.75
$ws.Range("N2").Value = -996.75
$ws.Range("H9").Value = 42.2
$ws.Range("I9").Value = 42.2
$ws.Range("K9").Value = 42.2
$ws.Range("M9").Value = 126.8
$ws.Range("H18").Value = 999
$ws.Range("I18").Value = 999
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 999
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -715
$ws.Range("N18").ClearContents()
$ws.Range("H33").Value = 597.94446
$ws.Range("I33").Value = 524.4666999999999
$ws.Range("J33").Value = 965.3333
$ws.Range("K33").Value = 524.4666999999999
$ws.Range("L33").Value = 965.3333
$ws.Range("M33").Value = -295.4666999999999
$ws.Range("N33").Value = -1423.3333
$ws.Range("H51").Value = 1000
$ws.Range("I51").Value = 1000
$ws.Range("K51").Value = 1000
$ws.Range("M51").Value = -516
$ws.Range("H100").Value = 3666
$ws.Range("I100").Value = 3699.5
$ws.Range("J100").Value = 3649.25
$ws.Range("K100").Value = 3699.5
$ws.Range("L100").Value = 3649.25
$ws.Range("M100").Value = -3158.5
$ws.Range("N100").Value = -4731.25
$ws.Range("H131").Value = 794
$ws.Range("I131").Value = 794
$ws.Range("K131").Value = 2382
$ws.Range("M131").Value = 2658
$ws.Range("H132").Value = 135415.6
$ws.Range("I132").Value = 135415.6
$ws.Range("K132").Value = 406246.8
$ws.Range("M132").Value = -403716.8

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 89017.8
$ws.Range("J92").Value = 89017.8
$ws.Range("L92").Value = 89017.8
$ws.Range("N92").Value = -94009.8

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 12964
$ws.Range("I97").Value = 12964
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 12964
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -11973
$ws.Range("N97").ClearContents()
$ws.Range("H102").Value = 12400
$ws.Range("I102").Value = 12400
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 12400
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -9155
$ws.Range("N102").ClearContents()

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1435.6666
$ws.Range("I5").Value = 807
$ws.Range("J5").Value = 1750
$ws.Range("K5").Value = 807
$ws.Range("L5").Value = 1750
$ws.Range("M5").Value = -695
$ws.Range("N5").Value = -1974
$ws.Range("H7").Value = 194.78572
$ws.Range("I7").Value = 174.88889
$ws.Range("K7").Value = 174.88889
$ws.Range("M7").Value = -61.88889
$ws.Range("H16").Value = 800
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H22").Value = 614.5833
$ws.Range("I22").Value = 698.3
$ws.Range("J22").Value = 196
$ws.Range("K22").Value = 698.3
$ws.Range("L22").Value = 196
$ws.Range("M22").Value = -348.3
$ws.Range("N22").Value = -896
$ws.Range("H99").Value = 2001397.8
$ws.Range("I99").Value = 1667996.6
$ws.Range("J99").Value = 2501499.5
$ws.Range("K99").Value = 1667996.6
$ws.Range("L99").Value = 2501499.5
$ws.Range("M99").Value = -1666498.6
$ws.Range("N99").Value = -2504495.5
$ws.Range("H107").Value = 521.3333
$ws.Range("I107").Value = 440.0909
$ws.Range("J107").Value = 744.75
$ws.Range("K107").Value = 440.0909
$ws.Range("L107").Value = 744.75
$ws.Range("M107").Value = 1479.9091
$ws.Range("N107").Value = -4584.75
$ws.Range("H113").Value = 800
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H126").Value = 2001397.8
$ws.Range("I126").Value = 1667996.6
$ws.Range("J126").Value = 2501499.5
$ws.Range("K126").Value = 5003989.800000001
$ws.Range("L126").Value = 7504498.5
$ws.Range("M126").Value = -5001519.800000001
$ws.Range("N126").Value = -7509438.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 308.2857
$ws.Range("J12").Value = 537
$ws.Range("L12").Value = 1611
$ws.Range("N12").Value = -1957
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H63").Value = 814
$ws.Range("J63").Value = 814
$ws.Range("L63").Value = 2442
$ws.Range("N63").Value = -3940
$ws.Range("H66").Value = 814
$ws.Range("J66").Value = 814
$ws.Range("L66").Value = 7326
$ws.Range("N66").Value = -14814
$ws.Range("H103").Value = 350
$ws.Range("I103").Value = 200
$ws.Range("J103").Value = 425
$ws.Range("K103").Value = 600
$ws.Range("L103").Value = 1275
$ws.Range("M103").Value = 279
$ws.Range("N103").Value = -3033
$ws.Range("H117").Value = 2557.75
$ws.Range("J117").Value = 10032
$ws.Range("L117").Value = 30096
$ws.Range("N117").Value = -36980
$ws.Range("H121").Value = 600.625
$ws.Range("I121").Value = 203.125
$ws.Range("J121").Value = 998.125
$ws.Range("K121").Value = 609.375
$ws.Range("L121").Value = 2994.375
$ws.Range("M121").Value = 700.625
$ws.Range("N121").Value = -5614.375
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 63500
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H97").Value = 549.25
$ws.Range("I97").Value = 415.83334
$ws.Range("J97").Value = 949.5
$ws.Range("K97").Value = 415.83334
$ws.Range("L97").Value = 949.5
$ws.Range("M97").Value = 80.16665999999998
$ws.Range("N97").Value = -1941.5
$ws.Range("H113").Value = 1670.6364
$ws.Range("I113").Value = 1499.1111
$ws.Range("J113").Value = 2442.5
$ws.Range("K113").Value = 1499.1111
$ws.Range("L113").Value = 2442.5
$ws.Range("M113").Value = 670.8888999999999
$ws.Range("N113").Value = -6782.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 5174.5
$ws.Range("I21").Value = 342
$ws.Range("K21").Value = 342
$ws.Range("M21").Value = -168
$ws.Range("H132").Value = 1715.2858
$ws.Range("I132").Value = 1239.4
$ws.Range("K132").Value = 3718.2
$ws.Range("M132").Value = -1188.2

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()
$ws.Range("H96").Value = 1856.125
$ws.Range("I96").Value = 1212.5
$ws.Range("J96").Value = 2499.75
$ws.Range("K96").Value = 1212.5
$ws.Range("L96").Value = 2499.75
$ws.Range("M96").Value = 160.5
$ws.Range("N96").Value = -5245.75
$ws.Range("H113").Value = 211.23077
$ws.Range("I113").Value = 140.22223
$ws.Range("J113").Value = 371
$ws.Range("K113").Value = 420.66669
$ws.Range("L113").Value = 1113
$ws.Range("M113").Value = 1749.33331
$ws.Range("N113").Value = -5453
$ws.Range("H136").Value = 2377
$ws.Range("I136").Value = 2377
$ws.Range("K136").Value = 7131
$ws.Range("M136").Value = -4581
